$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (BGR-packed OLE values, matching existing palette fills)
$blue   = 12611584   # RGB(0,112,192)  -> fill FF0070C0
$yellow = 65535       # RGB(255,255,0) -> fill FFFFFF00

# --- Write new, never-seen-before text first (controls shared-string order) ---
$ws.Range("M1").Value = "4x3 Keypad"
$ws.Range("N4").Value = "Expansion Board J10 (15)"
$ws.Range("N9").Value = "Expansion Board J10 (17)"
$ws.Range("N5").Value = "Expansion Board J11 (6)"
$ws.Range("N3").Value = "Expansion Board J11 (8)"
$ws.Range("N7").Value = "Expansion Board J11 (10)"
$ws.Range("N12").Value = "Row"
$ws.Range("N13").Value = "Column "

# --- Header row (reuses existing shared strings) ---
$ws.Range("M2").Value = "PIN #:"
$ws.Range("N2").Value = "Connection"
$ws.Range("O2").Value = "Wire Color"

# --- Data rows 3-9 ---
$ws.Range("M3").Value = 1
$ws.Range("O3").Value = "Blue"
$ws.Range("O3").Interior.Color = $blue

$ws.Range("M4").Value = 2
$ws.Range("O4").Value = "Yellow"
$ws.Range("O4").Interior.Color = $yellow

$ws.Range("M5").Value = 3
$ws.Range("O5").Value = "Blue"
$ws.Range("O5").Interior.Color = $blue

$ws.Range("M6").Value = 4
$ws.Range("N6").Value = "Expansion Board J10 (21)"
$ws.Range("O6").Value = "Yellow"
$ws.Range("O6").Interior.Color = $yellow

$ws.Range("M7").Value = 5
$ws.Range("O7").Value = "Blue"
$ws.Range("O7").Interior.Color = $blue

$ws.Range("M8").Value = 6
$ws.Range("N8").Value = "Expansion Board J10 (19)"
$ws.Range("O8").Value = "Yellow"
$ws.Range("O8").Interior.Color = $yellow

$ws.Range("M9").Value = 7
$ws.Range("O9").Value = "Yellow"
$ws.Range("O9").Interior.Color = $yellow

# --- Row/Column legend rows 12-13 (N/O only, no M) ---
$ws.Range("O12").Value = "Yellow"
$ws.Range("O12").Interior.Color = $yellow

$ws.Range("O13").Value = "Blue"
$ws.Range("O13").Interior.Color = $blue

# --- Column widths to fit the new content (best-fit style) ---
$ws.Columns.Item(10).ColumnWidth = 22.1666666667
$ws.Columns.Item(11).ColumnWidth = 9.6666666667
$ws.Columns.Item(13).ColumnWidth = 10.1666666667
$ws.Columns.Item(14).ColumnWidth = 22.1666666667
$ws.Columns.Item(15).ColumnWidth = 9.6666666667

# --- Selection moves to M19 after the edit ---
$ws.Range("M19").Select()
